# Update the NPCs worksheet:
#  - Update x_position (I) / y_position (J) values for the existing NPC rows
#  - Append two new NPC rows (Mrs.Piper / The Enchanted Snowman)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NPCs")

# --- Update x/y position values for existing rows (columns I=9, J=10) ---
$ws.Cells.Item(2, 9).Value = 384
$ws.Cells.Item(2, 10).Value = 1536

$ws.Cells.Item(3, 9).Value = 288
$ws.Cells.Item(3, 10).Value = 352

$ws.Cells.Item(4, 9).Value = 1136
$ws.Cells.Item(4, 10).Value = 512

$ws.Cells.Item(5, 9).Value = 976
$ws.Cells.Item(5, 10).Value = 1648

$ws.Cells.Item(6, 9).Value = 1040
$ws.Cells.Item(6, 10).Value = 592

$ws.Cells.Item(7, 9).Value = 1376
$ws.Cells.Item(7, 10).Value = 1392

$ws.Cells.Item(8, 9).Value = 1344
$ws.Cells.Item(8, 10).Value = 304

$ws.Cells.Item(9, 9).Value = 1888
$ws.Cells.Item(9, 10).Value = 1456

$ws.Cells.Item(10, 9).Value = 2016
$ws.Cells.Item(10, 10).Value = 1664

$ws.Cells.Item(11, 9).Value = 784
$ws.Cells.Item(11, 10).Value = 1040

$ws.Cells.Item(12, 9).Value = 688
$ws.Cells.Item(12, 10).Value = 512

$ws.Cells.Item(14, 9).Value = 560
$ws.Cells.Item(14, 10).Value = 1440

$ws.Cells.Item(15, 9).Value = 864
$ws.Cells.Item(15, 10).Value = 608

$ws.Cells.Item(16, 9).Value = 1072
$ws.Cells.Item(16, 10).Value = 1328

$ws.Cells.Item(17, 9).Value = 1696
$ws.Cells.Item(17, 10).Value = 496

$ws.Cells.Item(18, 9).Value = 416
$ws.Cells.Item(18, 10).Value = 208

$ws.Cells.Item(19, 9).Value = 1360
$ws.Cells.Item(19, 10).Value = 704

$ws.Cells.Item(20, 9).Value = 304
$ws.Cells.Item(20, 10).Value = 1664

$ws.Cells.Item(21, 9).Value = 592
$ws.Cells.Item(21, 10).Value = 1904

$ws.Cells.Item(22, 9).Value = 1200
$ws.Cells.Item(22, 10).Value = 720

$ws.Cells.Item(23, 9).Value = 1216
$ws.Cells.Item(23, 10).Value = 112

$ws.Cells.Item(24, 9).Value = 944
$ws.Cells.Item(24, 10).Value = 1952

$ws.Cells.Item(25, 9).Value = 1328
$ws.Cells.Item(25, 10).Value = 496

$ws.Cells.Item(26, 9).Value = 1008
$ws.Cells.Item(26, 10).Value = 1360

$ws.Cells.Item(27, 9).Value = 560
$ws.Cells.Item(27, 10).Value = 544

$ws.Cells.Item(28, 9).Value = 1008
$ws.Cells.Item(28, 10).Value = 880

$ws.Cells.Item(29, 9).Value = 1168
$ws.Cells.Item(29, 10).Value = 784

$ws.Cells.Item(30, 9).Value = 736
$ws.Cells.Item(30, 10).Value = 944

$ws.Cells.Item(31, 9).Value = 1280
$ws.Cells.Item(31, 10).Value = 2128

$ws.Cells.Item(32, 9).Value = 1760
$ws.Cells.Item(32, 10).Value = 400

$ws.Cells.Item(33, 9).Value = 1344
$ws.Cells.Item(33, 10).Value = 2032

$ws.Cells.Item(35, 9).Value = 2336
$ws.Cells.Item(35, 10).Value = 1984

$ws.Cells.Item(36, 9).Value = 160
$ws.Cells.Item(36, 10).Value = 2192

# --- Append two new NPC rows ---
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "Mrs.Piper"
$ws.Cells.Item(37, 3).Value = "Mrs. Piper"
$ws.Cells.Item(37, 4).Value = 2
$ws.Cells.Item(37, 5).Value = "The Ice Plane"
$ws.Cells.Item(37, 9).Value = 640
$ws.Cells.Item(37, 10).Value = 1552

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "TheEnchantedSnowman"
$ws.Cells.Item(38, 3).Value = "The Enchanted Snowman"
$ws.Cells.Item(38, 4).Value = 2
$ws.Cells.Item(38, 5).Value = "The Ice Plane"
$ws.Cells.Item(38, 9).Value = 304
$ws.Cells.Item(38, 10).Value = 1232
